# Excel COM-interop edit script
# Implements: add two new worksheets "DualAxis_Columns" and "DualAxis_Grouped"
# at the end of the workbook (v0.7.5 "Second Y-axis support"), populate them
# with data, and make the last new sheet the active tab (mirrors removing
# tabSelected from the previously-active "xy_log" sheet).

$wb = $excel.ActiveWorkbook

# --- Add "DualAxis_Columns" as the new last sheet --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dualAxisColumns = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$dualAxisColumns.Name = "DualAxis_Columns"

$dualAxisColumnsData = @(
  @("Sample ID","Viability (%)","Protein (mg/mL)","Absorbance (AU)"),
  @(0,98.417727190009515,2.3171366678656748,[double]"5.1656317157017823E-2"),
  @(10,91.08065975224045,2.6030854868705391,0.13487772563561179),
  @(20,89.343627902203139,2.47108795339118,0.20760412881077361),
  @(30,89.393722471738243,3.5972251908184378,0.2890717051166809),
  @(40,87.476933159990821,3.8844753158751049,0.36446832512996991),
  @(50,83.34273656237994,4.2969233999082128,0.44401896687959658),
  @(60,79.768703435223514,5.1835028866522608,0.534062629111971),
  @(70,76.397792608821419,5.809299856748785,0.6167812001428542),
  @(80,71.042956019265148,6.2793840357348598,0.68963994939209838),
  @(90,69.560311583210577,6.2482347430332084,0.77501766448946019),
  @(100,67.078722458080421,6.9072362872446353,0.85180818012523829)
)

$r = 1
foreach ($row in $dualAxisColumnsData) {
    $c = 1
    foreach ($val in $row) {
        $dualAxisColumns.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
$dualAxisColumns.Range("A1:D1048576").Select() | Out-Null

# --- Add "DualAxis_Grouped" as the new last sheet (becomes active tab) -----
$dualAxisGrouped = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dualAxisColumns)
$dualAxisGrouped.Name = "DualAxis_Grouped"

$dualAxisGroupedData = @(
  @("Dose (uM)","Response","Cell Line"),
  @(0.1,533.90919744193639,"HEK293"),
  @(0.1,0.54706009933035715,"CHO"),
  @(0.1,541.45791669640937,"HEK293"),
  @(0.1,0.87931893084619928,"CHO"),
  @(0.1,497.3658210552079,"HEK293"),
  @(0.1,0.52788556817475796,"CHO"),
  @(0.3,624.94078582294287,"HEK293"),
  @(0.3,0.92056618979303595,"CHO"),
  @(0.3,522.50642575082566,"HEK293"),
  @(0.3,0.8755912526796471,"CHO"),
  @(0.3,522.80926040695022,"HEK293"),
  @(0.3,0.67393329324840279,"CHO"),
  @(1,754.1838302729152,"HEK293"),
  @(1,1.0672736187182561,"CHO"),
  @(1,655.83982506896177,"HEK293"),
  @(1,1.337472161801621,"CHO"),
  @(1,691.44416067789223,"HEK293"),
  @(1,1.512779134168871,"CHO"),
  @(3,1204.598796223939,"HEK293"),
  @(3,2.4675392597329422,"CHO"),
  @(3,1323.282438446078,"HEK293"),
  @(3,2.7048447399026929,"CHO"),
  @(3,1253.376410234396,"HEK293"),
  @(3,2.4650503627573088,"CHO"),
  @(10,1761.090061227673,"HEK293"),
  @(10,4.1635477071850646,"CHO"),
  @(10,1730.7595185828161,"HEK293"),
  @(10,4.2165027929122258,"CHO"),
  @(10,1758.2772629579911,"HEK293"),
  @(10,4.0830244392844364,"CHO"),
  @(30,1923.9345243409571,"HEK293"),
  @(30,5.1033763929393334,"CHO"),
  @(30,1953.3449937155301,"HEK293"),
  @(30,4.5213785702463651,"CHO"),
  @(30,1995.1471005575861,"HEK293"),
  @(30,4.4887520260433407,"CHO")
)

$r = 1
foreach ($row in $dualAxisGroupedData) {
    $c = 1
    foreach ($val in $row) {
        $dualAxisGrouped.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Make the newly added sheet the active/selected sheet (also clears
# tabSelected on whichever sheet was previously active, e.g. "xy_log").
$dualAxisGrouped.Activate() | Out-Null
$dualAxisGrouped.Range("U53").Select() | Out-Null
